$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 221
$ws1.Range("F5").Value = 1011
$ws1.Range("F6").Value = 5650
$ws1.Range("F7").Value = 520
$ws1.Range("F8").Value = 731
$ws1.Range("F9").Value = 973
$ws1.Range("F15").Value = 25
$ws1.Range("F17").Value = 1905
$ws1.Range("F18").Value = 1487
$ws1.Range("F19").Value = 970
$ws1.Range("F22").Value = 350
$ws1.Range("F23").Value = 575
$ws1.Range("F24").Value = 168
$ws1.Range("F25").Value = 1061
$ws1.Range("F27").Value = 526
$ws1.Range("F28").Value = 3107
$ws1.Range("F30").Value = 112
$ws1.Range("F31").Value = 71
$ws1.Range("F32").Value = 135
$ws1.Range("F34").Value = 424
$ws1.Range("F39").Value = 304
$ws1.Range("F40").Value = 754
$ws1.Range("F42").Value = 60
$ws1.Range("F43").Value = 64
$ws1.Range("F44").Value = 72

# Sheet 2: 演出 (Performance) - column F updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 218
$ws2.Range("F6").Value = 149

# Sheet 4: 全部类型 (All Types) - column F updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 221
$ws4.Range("F5").Value = 1011
$ws4.Range("F7").Value = 5650
$ws4.Range("F8").Value = 520
$ws4.Range("F9").Value = 731
$ws4.Range("F11").Value = 218
$ws4.Range("F12").Value = 973
$ws4.Range("F15").Value = 149
$ws4.Range("F20").Value = 25
$ws4.Range("F23").Value = 1905
$ws4.Range("F24").Value = 1487
$ws4.Range("F25").Value = 970
$ws4.Range("F27").Value = 350
$ws4.Range("F29").Value = 575
$ws4.Range("F30").Value = 168
$ws4.Range("F31").Value = 1061
$ws4.Range("F32").Value = 3107
$ws4.Range("F34").Value = 112
$ws4.Range("F35").Value = 71
$ws4.Range("F36").Value = 135
$ws4.Range("F38").Value = 424
$ws4.Range("F42").Value = 304
$ws4.Range("F43").Value = 754
$ws4.Range("F45").Value = 64
$ws4.Range("F46").Value = 72
